$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-slot labels in column C (the "Время (мск)" column).
# Rows 4 and 5 ("10:55-11:0" / "11:0-11:5") are left untouched.
$ws.Range("C2").Value = "2:55-3:0"
$ws.Range("C3").Value = "3:0-3:5"
$ws.Range("C6").Value = "18:55-19:0"
$ws.Range("C7").Value = "19:0-19:5"

# Move the active selection from C11 to B11, as recorded in the saved file.
$ws.Range("B11").Select()
